$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core numeric change: D13 (master of a long dependent-formula chain: D14=D13,
# D15..D25 shared si=1 chain back to D14, D31=D13, D32..D35 shared si=2 chain back
# to D14) goes from 4 to 8. Every J/K value down the sheet recalculates from this.
$ws.Range("D13").Value = 8

# --- Row 35: B35/C35 16 -> 24 (D35 already cascades from D13 above)
$ws.Range("B35").Value = 24
$ws.Range("C35").Value = 24

# --- Row 37: building-cost row tweaks
$ws.Range("C37").Value = 24
$ws.Range("E37").Value = 16

# --- Remove the "Human building" extras no longer tracked (O7, P7, Q7, O12)
$ws.Range("O7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("O12").ClearContents()

# --- Delete the now-unused building rows 74:82 (Build basic/advanced structure,
# Human Farm/Hall/Barracks/Mill/Smith/Church/Stable) - also drops their shared
# strings on save.
$ws.Rows("74:82").Delete()

# --- Formulas that must be retyped explicitly (Excel doesn't auto-shrink these
# ranges just because rows were deleted elsewhere / summation range changed)
$ws.Range("H37").Formula = "=COUNTA(A38:A80)"
$ws.Range("P35").Formula = "=(SUM(K5:K187)+Q4+Q5+Q6+S10)/1024"

# --- View state: restore selection to B6 (matches author's saved cursor position)
$ws.Range("B6").Select()
